$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the two date strings (shared across their respective rows)
$ws.Range("B2:B10").Value = "02/44/2025"
$ws.Range("B11:B12").Value = "02/45/2025"

# Update Monto (amount) values in rows 2-10
$ws.Range("D2").Value = 2250.0
$ws.Range("D4").Value = 100.0
$ws.Range("D5").Value = 1780.0
$ws.Range("D6").Value = 80.0
$ws.Range("D7").Value = 1200.0
$ws.Range("D9").Value = 10000.0
$ws.Range("D10").Value = 1300.0

# Row 11: id becomes 10, product becomes Chocolate, amount becomes 3560
$ws.Range("A11").Value = 10.0
$ws.Range("C11").Value = "Chocolate"
$ws.Range("D11").Value = 3560.0

# Row 12: id becomes 11, product becomes Alfajor, amount becomes 1300
$ws.Range("A12").Value = 11.0
$ws.Range("C12").Value = "Alfajor"
$ws.Range("D12").Value = 1300.0

# Remove the now-obsolete trailing rows 13-19
$ws.Range("A13:D19").EntireRow.Delete()
